# Update cryptocurrency price/volume table per latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.686.14"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.919.59"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'239.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.4940"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.2978"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "'0.06773"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").Value = "1.887.87"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'17.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "'0.07346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "'5.168"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").Value = "'88.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").Value = "'0.6727"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").Value = "30.664.94"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'0.000007954"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "'13.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.12%  "
$ws.Range("D20").Value = "2.153.68"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "'5.349"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.99%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'199.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.30%  "
$ws.Range("D24").Value = "'6.310"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.74%  "
$ws.Range("D25").Value = "'9.657"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("D26").Value = "'166.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.75%  "
$ws.Range("D27").Value = "'18.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").Value = "'1.967"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("D29").Value = "'1.482"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.97%  "
$ws.Range("D30").Value = "'4.370"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "'0.09188"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").Value = "'4.065"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").Value = "'0.05290"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "'0.7441"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "'2.732"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "'0.01840"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'2.724"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("D39").Value = "'0.9249"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("D40").Value = "'2.084"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("D41").Value = "'74.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +29.35%  "
$ws.Range("D42").Value = "'0.4460"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").Value = "'5.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.04%  "
$ws.Range("D44").Value = "'106.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").Value = "'1.003"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").Value = "'0.1389"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("D47").Value = "'7.628"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'36.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.16%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.050"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.80%  "
$ws.Range("D50").Value = "'0.05882"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'0.4046"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.67%  "
